$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "GPU mining additionals" helper row (old C10:D10) up into I8:J8 ---
# Copy formatting + value from C10:D10 into I8:J8, then restore the formula in J8
# (Copy alone drops the formula, so we set it explicitly afterwards).
$ws.Range("C10:D10").Copy($ws.Range("I8:J8"))
$ws.Cells.Item(8, 10).Formula = "=9000/15"

# --- Move the other helper formula (old D11) up into J9 ---
$ws.Range("D11").Copy($ws.Range("J9"))
$ws.Cells.Item(9, 10).Formula = "=450+29+75+648+120"

# --- Remove the now-empty trailing rows 10 and 11 ---
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()

# --- Give column I (the new helper label column) an explicit width ---
$ws.Columns.Item(9).ColumnWidth = 22.75

# --- Update the selection to match the author's final cursor position ---
$ws.Range("J11").Select()
